$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-11-30"

# Update the label for November in column A (row 12)
$ws.Range("A12").Value = "November (through 11-30)"

# Update November row (row 12) values
$ws.Range("B12").Value = 33
$ws.Range("C12").Value = 77
$ws.Range("D12").Value = 111
$ws.Range("F12").Value = 52
$ws.Range("G12").Value = 207
$ws.Range("H12").Value = 202

# Update Total row (row 13) values
$ws.Range("B13").Value = 291
$ws.Range("C13").Value = 563
$ws.Range("D13").Value = 821
$ws.Range("F13").Value = 534
$ws.Range("G13").Value = 1264
$ws.Range("H13").Value = 1645
